# Apply crypto price/volume update (commit: Updated cryptos list on Sat Jan  6 15:45:58 UTC 2024 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "95.50" or "6.37"
# are not silently coerced to numbers by Excel (matches original inlineStr cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "44.139.77"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "2.244.74"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "306.76"
$ws.Range("E5").Value = "  -2.06%  "
$ws.Range("D6").Value = "95.50"
$ws.Range("E6").Value = "  -3.23%  "
$ws.Range("E7").Value = "  +0.93%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("D10").Value = "34.92"
$ws.Range("E10").Value = "  -3.53%  "
$ws.Range("D11").Value = "0.0816"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "7.24"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.588.01"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.323.37"
$ws.Range("E15").Value = "  +3.83%  "
$ws.Range("D16").Value = "0.832"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").Value = "13.61"
$ws.Range("E17").Value = "  -3.36%  "
$ws.Range("D18").Value = "44.071.46"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "6.37"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "12.14"
$ws.Range("E21").Value = "  -6.30%  "
$ws.Range("D22").Value = "65.43"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "236.51"
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "9.92"
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("D29").Value = "37.36"
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").Value = "5.99"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").Value = "20.05"
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("D32").Value = "152.24"
$ws.Range("E32").Value = "  -3.58%  "
$ws.Range("D33").Value = "0.0804"
$ws.Range("E33").Value = "  -3.10%  "
$ws.Range("E34").Value = "  +3.92%  "
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("D38").Value = "1.76"
$ws.Range("E38").Value = "  -6.90%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "3.41"
$ws.Range("E39").Value = "  -5.47%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "3.87"
$ws.Range("E40").Value = "  -4.04%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "14.46"
$ws.Range("E41").Value = "  -7.92%  "
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "1.736.84"
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("D45").Value = "82.83"
$ws.Range("E45").Value = "  +3.45%  "
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("D47").Value = "100.04"
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("E48").Value = "  -4.57%  "
$ws.Range("D49").Value = "8.12"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").Value = "54.66"
$ws.Range("E50").Value = "  -3.06%  "
$ws.Range("D51").Value = "67.87"
$ws.Range("E51").Value = "  -7.20%  "

# Restore the default style on column D so no stray number-format styling is left behind.
$ws.Range("D2:D51").Style = "Normal"
